$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update renamed location labels in column A ---
$ws.Range("A8").Value = "大食堂"
$ws.Range("A18").Value = "纳米楼"
$ws.Range("A20").Value = "医学院药学院"
$ws.Range("A21").Value = "金工化学生物"

# --- Add header labels for the new average columns ---
$ws.Range("F4").Value = "经度"
$ws.Range("K4").Value = "纬度"

# --- Fill in the AVERAGE formulas down columns F (longitude avg) and K (latitude avg) ---
for ($r = 5; $r -le 25; $r++) {
    $ws.Range("F$r").Formula = "=AVERAGE(B${r}:E${r})"
    $ws.Range("K$r").Formula = "=AVERAGE(G${r}:J${r})"
}

# --- Update the active selection to match the saved view state ---
$ws.Range("K21").Select() | Out-Null
